# The sheet originally has 8 columns (A-H) with data in rows 1-2.
# The desired result keeps only the former column F ("input_Name" / empty),
# which becomes the sole column A, with column width 12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns A-E (div_backdropElements_class ... div_integrationIcons_internalRoleTreeitemName)
# so that former column F ("input_Name") becomes column A.
$ws.Range("A1:E1").EntireColumn.Delete()

# Former columns G:H (p_sessionInfo_class, p_sessionInfo_class_1) are now at B:C; remove them.
$ws.Range("B1:C1").EntireColumn.Delete()

# Match the target column width of 12 characters (Excel applies a small internal
# padding offset when setting ColumnWidth, so compensate to land exactly on 12).
$ws.Columns("A").ColumnWidth = 11.166666666666666

Write-Host ("UsedRange: " + $ws.UsedRange.Address())
Write-Host ("A1: " + $ws.Range("A1").Text)
Write-Host ("A2: " + $ws.Range("A2").Text)
